# Apply updated odds values to row 3 of the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.7
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 2.38
$ws.Range("L3").Value = 5
$ws.Range("X3").Value = 8
$ws.Range("Z3").Value = 13
$ws.Range("AB3").Value = 29
$ws.Range("AC3").Value = 9.5
$ws.Range("AE3").Value = 17
$ws.Range("AG3").Value = 301
$ws.Range("AH3").Value = 13
$ws.Range("AI3").Value = 26
$ws.Range("AJ3").Value = 17
$ws.Range("AL3").Value = 41
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 9
$ws.Range("AQ3").Value = 29
$ws.Range("AW3").Value = 6.5
$ws.Range("AX3").Value = 26
$ws.Range("AZ3").Value = 101
$ws.Range("BA3").Value = 126
